$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2000656.6
$ws.Range("I38").Value = 820.75
$ws.Range("J38").Value = 10000000
$ws.Range("K38").Value = 2462.25
$ws.Range("L38").Value = 30000000
$ws.Range("M38").Value = -2090.25
$ws.Range("N38").Value = -30000744

$ws.Range("H47").Value = 2000
$ws.Range("I47").Value = 2000
$ws.Range("K47").Value = 2000
$ws.Range("M47").Value = -1028

$ws.Range("H64").Value = 2266052
$ws.Range("I64").Value = 3348038.5
$ws.Range("J64").Value = 3716.6365
$ws.Range("K64").Value = 3348038.5
$ws.Range("L64").Value = 3716.6365
$ws.Range("M64").Value = -3347790.5
$ws.Range("N64").Value = -4212.636500000001

$ws.Range("H67").Value = 2266052
$ws.Range("I67").Value = 3348038.5
$ws.Range("J67").Value = 3716.6365
$ws.Range("K67").Value = 3348038.5
$ws.Range("L67").Value = 3716.6365
$ws.Range("M67").Value = -3347180.5
$ws.Range("N67").Value = -5432.636500000001

$ws.Range("H74").Value = 4421.3687
$ws.Range("I74").Value = 3778
$ws.Range("J74").Value = 5000.4
$ws.Range("K74").Value = 3778
$ws.Range("L74").Value = 5000.4
$ws.Range("M74").Value = -2842
$ws.Range("N74").Value = -6872.4

$ws.Range("H76").Value = 3491.6487
$ws.Range("I76").Value = 3457.6072
$ws.Range("J76").Value = 3597.5557
$ws.Range("K76").Value = 3457.6072
$ws.Range("L76").Value = 3597.5557
$ws.Range("M76").Value = -3142.6072
$ws.Range("N76").Value = -4227.5557

$ws.Range("H77").Value = 4421.3687
$ws.Range("I77").Value = 3778
$ws.Range("J77").Value = 5000.4
$ws.Range("K77").Value = 18890
$ws.Range("L77").Value = 25002
$ws.Range("M77").Value = -14210
$ws.Range("N77").Value = -34362

$ws.Range("H79").Value = 3491.6487
$ws.Range("I79").Value = 3457.6072
$ws.Range("J79").Value = 3597.5557
$ws.Range("K79").Value = 3457.6072
$ws.Range("L79").Value = 3597.5557
$ws.Range("M79").Value = -2365.6072
$ws.Range("N79").Value = -5781.5557

$ws.Range("H100").Value = 2342.25
$ws.Range("I100").Value = 1645.8
$ws.Range("J100").Value = 3503
$ws.Range("K100").Value = 1645.8
$ws.Range("L100").Value = 3503
$ws.Range("M100").Value = -1104.8
$ws.Range("N100").Value = -4585

$ws.Range("H138").Value = 2066.8667
$ws.Range("J138").Value = 2999.1396
$ws.Range("L138").Value = 8997.418799999999
$ws.Range("N138").Value = -19277.4188

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5773.3916
$ws.Range("I32").Value = 3783.753
$ws.Range("K32").Value = 3783.753
$ws.Range("M32").Value = -3496.753

$ws.Range("H45").Value = 1809.6923
$ws.Range("I45").Value = 1921.5
$ws.Range("J45").Value = 1630.8
$ws.Range("K45").Value = 1921.5
$ws.Range("L45").Value = 1630.8
$ws.Range("M45").Value = -1544.5
$ws.Range("N45").Value = -2384.8

$ws.Range("H106").Value = 47713.332
$ws.Range("J106").Value = 47713.332
$ws.Range("L106").Value = 47713.332
$ws.Range("N106").Value = -50237.332

$ws.Range("H117").Value = 32432.334
$ws.Range("J117").Value = 32432.334
$ws.Range("L117").Value = 32432.334
$ws.Range("N117").Value = -41610.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2667.375
$ws.Range("I62").Value = 2186.6667
$ws.Range("J62").Value = 3285.4285
$ws.Range("K62").Value = 2186.6667
$ws.Range("L62").Value = 3285.4285
$ws.Range("M62").Value = -1562.6667
$ws.Range("N62").Value = -4533.4285

$ws.Range("H65").Value = 2667.375
$ws.Range("I65").Value = 2186.6667
$ws.Range("J65").Value = 3285.4285
$ws.Range("K65").Value = 10933.3335
$ws.Range("L65").Value = 16427.1425
$ws.Range("M65").Value = -7813.333500000001
$ws.Range("N65").Value = -22667.1425

$ws.Range("H107").Value = 476.53125
$ws.Range("I107").Value = 284.31818
$ws.Range("J107").Value = 899.4
$ws.Range("K107").Value = 284.31818
$ws.Range("L107").Value = 899.4
$ws.Range("M107").Value = 1635.68182
$ws.Range("N107").Value = -4739.4

$ws.Range("H132").Value = 12822851
$ws.Range("I132").Value = 16668476
$ws.Range("K132").Value = 50005428
$ws.Range("M132").Value = -50002898

$ws.Range("H140").Value = 39237.9
$ws.Range("J140").Value = 39237.9
$ws.Range("L140").Value = 39237.9
$ws.Range("N140").Value = -49597.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 11000
$ws.Range("J105").Value = 11000
$ws.Range("L105").Value = 33000
$ws.Range("N105").Value = -38242

$ws.Range("H113").Value = 40000800
$ws.Range("I113").Value = 250000500
$ws.Range("J113").Value = 856.619
$ws.Range("K113").Value = 750001500
$ws.Range("L113").Value = 2569.857
$ws.Range("M113").Value = -749999330
$ws.Range("N113").Value = -6909.857

$ws.Range("H129").Value = 3910.2222
$ws.Range("I129").Value = 3898.7693
$ws.Range("J129").Value = 3920.8572
$ws.Range("K129").Value = 11696.3079
$ws.Range("L129").Value = 11762.5716
$ws.Range("M129").Value = -6696.3079
$ws.Range("N129").Value = -21762.5716

$ws.Range("H131").Value = 1830.3651
$ws.Range("I131").Value = 2599.4211
$ws.Range("J131").Value = 1498.2727
$ws.Range("K131").Value = 7798.263300000001
$ws.Range("L131").Value = 4494.8181
$ws.Range("M131").Value = -2758.263300000001
$ws.Range("N131").Value = -14574.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11537.346
$ws.Range("I70").Value = 14887.277
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 14887.277
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -14617.277
$ws.Range("N70").Value = -4540

$ws.Range("H73").Value = 11537.346
$ws.Range("I73").Value = 14887.277
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 14887.277
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -13951.277
$ws.Range("N73").Value = -5872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3693.9119
$ws.Range("I40").Value = 3635
$ws.Range("J40").Value = 3778.0715
$ws.Range("K40").Value = 3635
$ws.Range("L40").Value = 3778.0715
$ws.Range("M40").Value = -3499
$ws.Range("N40").Value = -4050.0715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1450
$ws.Range("I81").Value = 1900
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 3800
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -2739
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 1450
$ws.Range("I84").Value = 1900
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 19000
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -13696
$ws.Range("N84").Value = -20608

$ws.Range("H132").Value = 1746.8928
$ws.Range("I132").Value = 1445.5
$ws.Range("J132").Value = 3555.25
$ws.Range("K132").Value = 4336.5
$ws.Range("L132").Value = 10665.75
$ws.Range("M132").Value = -1806.5
$ws.Range("N132").Value = -15725.75
